# Update computed price/profit figures on the per-craft-class Leve sheets.
# Mirrors a scheduled market-data refresh run (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 525.04877
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 525.04877
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1575.14631
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1911.14631

# Row 39
$ws.Range("H39").Value = 1051.1538
$ws.Range("I39").Value = 170.16667
$ws.Range("J39").Value = 1806.2858
$ws.Range("K39").Value = 510.50001
$ws.Range("L39").Value = 5418.857400000001
$ws.Range("M39").Value = -214.50001
$ws.Range("N39").Value = -6010.857400000001

# Row 137
$ws.Range("H137").Value = 24432438
$ws.Range("I137").Value = 55558372
$ws.Range("K137").Value = 166675116
$ws.Range("M137").Value = -166672566

# Row 138
$ws.Range("H138").Value = 7685.12
$ws.Range("I138").Value = 4158
$ws.Range("J138").Value = 9056.777
$ws.Range("K138").Value = 12474
$ws.Range("L138").Value = 27170.331
$ws.Range("M138").Value = -7334
$ws.Range("N138").Value = -37450.331

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1871.7755
$ws.Range("I32").Value = 1660.3778
$ws.Range("K32").Value = 1660.3778
$ws.Range("M32").Value = -1373.3778

# Row 43
$ws.Range("H43").Value = 38564.168
$ws.Range("J43").Value = 38248.75
$ws.Range("L43").Value = 38248.75
$ws.Range("N43").Value = -38874.75

# Row 45
$ws.Range("H45").Value = 2799.75
$ws.Range("I45").Value = 1399.6666
$ws.Range("K45").Value = 1399.6666
$ws.Range("M45").Value = -1022.6666

# Row 60
$ws.Range("H60").Value = 71278.79
$ws.Range("I60").Value = 71278.79
$ws.Range("K60").Value = 71278.79
$ws.Range("M60").Value = -70545.79

# Row 61
$ws.Range("H61").Value = 78890170
$ws.Range("I61").Value = 116667920
$ws.Range("K61").Value = 116667920
$ws.Range("M61").Value = -116667708

# Row 136
$ws.Range("H136").Value = 78890170
$ws.Range("I136").Value = 116667920
$ws.Range("K136").Value = 350003760
$ws.Range("M136").Value = -350001210

$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Range("H6").Value = 61666
$ws.Range("J6").Value = 61666
$ws.Range("L6").Value = 61666
$ws.Range("N6").Value = -61892

# Row 105
$ws.Range("H105").Value = 563661.06
$ws.Range("I105").Value = 759182.7
$ws.Range("J105").Value = 9683
$ws.Range("K105").Value = 759182.7
$ws.Range("L105").Value = 9683
$ws.Range("M105").Value = -757435.7
$ws.Range("N105").Value = -13177

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 37414500
$ws.Range("I31").Value = 45457764
$ws.Range("K31").Value = 45457764
$ws.Range("M31").Value = -45457469

# Row 34
$ws.Range("H34").Value = 37414500
$ws.Range("I34").Value = 45457764
$ws.Range("K34").Value = 45457764
$ws.Range("M34").Value = -45457562

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 87
$ws.Range("H87").Value = 99000
$ws.Range("J87").Value = 99000
$ws.Range("L87").Value = 99000
$ws.Range("N87").Value = -101372

# Row 90
$ws.Range("H90").Value = 99000
$ws.Range("J90").Value = 99000
$ws.Range("L90").Value = 297000
$ws.Range("N90").Value = -308856

# Row 107
$ws.Range("H107").Value = 1502.0817
$ws.Range("I107").Value = 1311.2368
$ws.Range("J107").Value = 2161.3635
$ws.Range("K107").Value = 1311.2368
$ws.Range("L107").Value = 2161.3635
$ws.Range("M107").Value = 608.7632000000001
$ws.Range("N107").Value = -6001.363499999999

# Row 111
$ws.Range("H111").Value = 64999
$ws.Range("J111").Value = 64999
$ws.Range("L111").Value = 64999
$ws.Range("N111").Value = -73179

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3002
$ws.Range("J5").Value = 4999
$ws.Range("L5").Value = 14997
$ws.Range("N5").Value = -15221

# Row 34
$ws.Range("H34").Value = 10867.833
$ws.Range("I34").Value = 299
$ws.Range("J34").Value = 16152.25
$ws.Range("K34").Value = 897
$ws.Range("L34").Value = 48456.75
$ws.Range("M34").Value = -813
$ws.Range("N34").Value = -48624.75

# Row 44
$ws.Range("H44").Value = 5087.778
$ws.Range("I44").Value = 173.16667
$ws.Range("J44").Value = 14917
$ws.Range("K44").Value = 519.50001
$ws.Range("L44").Value = 44751
$ws.Range("M44").Value = -121.50001
$ws.Range("N44").Value = -45547

# Row 86
$ws.Range("H86").Value = 1273
$ws.Range("J86").Value = 1273
$ws.Range("L86").Value = 3819
$ws.Range("N86").Value = -6191

# Row 89
$ws.Range("H89").Value = 1273
$ws.Range("J89").Value = 1273
$ws.Range("L89").Value = 11457
$ws.Range("N89").Value = -23313

# Row 135
$ws.Range("H135").Value = 3002
$ws.Range("J135").Value = 4999
$ws.Range("L135").Value = 44991
$ws.Range("N135").Value = -50061

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 61.11111
$ws.Range("I2").Value = 69.166664
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 69.166664
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 43.833336
$ws.Range("N2").Value = -271

# Row 113
$ws.Range("H113").Value = 843862.06
$ws.Range("I113").Value = 2177
$ws.Range("K113").Value = 2177
$ws.Range("M113").Value = -7

# Row 129
$ws.Range("H129").Value = 99999
$ws.Range("J129").Value = 99999
$ws.Range("L129").Value = 99999
$ws.Range("N129").Value = -109999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5510.8
$ws.Range("I40").Value = 5388.5
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 5388.5
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -5252.5
$ws.Range("N40").Value = -6272

# Row 68
$ws.Range("H68").Value = 8103968.5
$ws.Range("I68").Value = 10804224
$ws.Range("J68").Value = 3200.6667
$ws.Range("K68").Value = 10804224
$ws.Range("L68").Value = 3200.6667
$ws.Range("M68").Value = -10803475
$ws.Range("N68").Value = -4698.6667

# Row 71
$ws.Range("H71").Value = 8103968.5
$ws.Range("I71").Value = 10804224
$ws.Range("J71").Value = 3200.6667
$ws.Range("K71").Value = 54021120
$ws.Range("L71").Value = 16003.3335
$ws.Range("M71").Value = -54017376
$ws.Range("N71").Value = -23491.3335

# Row 93
$ws.Range("H93").Value = 2419934.2
$ws.Range("I93").Value = 2183.3125
$ws.Range("J93").Value = 7946222
$ws.Range("K93").Value = 2183.3125
$ws.Range("L93").Value = 7946222
$ws.Range("M93").Value = -935.3125
$ws.Range("N93").Value = -7948718

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 40495
$ws.Range("I51").Value = 70
$ws.Range("J51").Value = 48580
$ws.Range("K51").Value = 70
$ws.Range("L51").Value = 48580
$ws.Range("M51").Value = 440
$ws.Range("N51").Value = -49600

# Row 56
$ws.Range("H56").Value = 65000
$ws.Range("J56").Value = 65000
$ws.Range("L56").Value = 65000
$ws.Range("N56").Value = -66428
